$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.0498220640569395
$wsSummary.Range("C2").Value = 0.0498220640569395
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.09491525423728814
$wsSummary.Range("F2").Value = 0.2077151335311573
$wsSummary.Range("G2").Value = 0.5768621236133122
$wsSummary.Range("H2").Value = 0.7486623863028358
$wsSummary.Range("I2").Value = 28
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$wsClass = $wb.Worksheets.Item("Classification Report")

$wsClass.Range("B2").Value = 0
$wsClass.Range("C2").Value = 0
$wsClass.Range("D2").Value = 0

$wsClass.Range("B3").Value = 0.0498220640569395
$wsClass.Range("C3").Value = 1
$wsClass.Range("D3").Value = 0.09491525423728814

$wsClass.Range("B4").Value = 0.0498220640569395
$wsClass.Range("C4").Value = 0.0498220640569395
$wsClass.Range("D4").Value = 0.0498220640569395
$wsClass.Range("E4").Value = 0.0498220640569395

$wsClass.Range("B5").Value = 0.02491103202846975
$wsClass.Range("C5").Value = 0.5
$wsClass.Range("D5").Value = 0.04745762711864407

$wsClass.Range("B6").Value = 0.002482238066893783
$wsClass.Range("C6").Value = 0.0498220640569395
$wsClass.Range("D6").Value = 0.004728873876590867

# --- Sheet: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

$wsConf.Range("B2").Value = 0
$wsConf.Range("C2").Value = 534

$wsConf.Range("B3").Value = 0
$wsConf.Range("C3").Value = 28
